$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 27380
$ws.Range("B3").Value = 1880
$ws.Range("B4").Value = 15535
$ws.Range("B5").Value = 17569
